# Event_Sorting.xlsx update:
#  - Move the 4 previously "orphaned" SRM status events (rows 60, 62, 63, 66,
#    which had no running-total value in column B) up into the main sorted
#    list, right after the existing header/seed rows (A1:B3), i.e. before the
#    former A4.
#  - The running total formula in column B (B*=prev+100) is extended so it
#    now also covers these 4 newly (re)integrated rows, and the whole list
#    grows from A4:B39 to A4:B43.
#  - Adjust sheet view selection & dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture the text of the 4 orphan rows before touching anything.
# ---------------------------------------------------------------------
$orphan1 = $ws.Range("A60").Value2   # SRM: Change was Transmitted
$orphan2 = $ws.Range("A62").Value2   # SRM: Held
$orphan3 = $ws.Range("A63").Value2   # SRM: Deleted
$orphan4 = $ws.Range("A66").Value2   # SRM: Awaiting Approval

# ---------------------------------------------------------------------
# 2. Remove the old orphan rows (bottom-up so row numbers stay valid).
# ---------------------------------------------------------------------
$ws.Rows("66:66").Delete()
$ws.Rows("63:63").Delete()
$ws.Rows("62:62").Delete()
$ws.Rows("60:60").Delete()

# ---------------------------------------------------------------------
# 3. Insert 4 fresh rows right before the old row 4, pushing the sorted
#    list (formerly A4:B39) down to A8:B43.
# ---------------------------------------------------------------------
$ws.Rows("4:7").Insert()

# ---------------------------------------------------------------------
# 4. Populate the newly inserted rows with the orphan event names and a
#    running-total formula consistent with the rest of the column.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = $orphan1
$ws.Range("A5").Value = $orphan2
$ws.Range("A6").Value = $orphan3
$ws.Range("A7").Value = $orphan4

$ws.Range("B4").Formula = "=B3+100"
$ws.Range("B5").Formula = "=B4+100"
$ws.Range("B6").Formula = "=B5+100"
$ws.Range("B7").Formula = "=B6+100"

# The row that used to be the anchor of the shared formula (old row 4,
# now shifted to row 8) keeps its original formula text after the
# insert, so make sure it still correctly refers to the row above it.
# Re-apply the same "=prev+100" pattern all the way down to the new
# last row (43) so the whole B4:B43 range is recognised as one
# consistent (and re-sharable) formula run.
for ($r = 8; $r -le 43; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev+100"
}

# ---------------------------------------------------------------------
# 5. Update the sheet view: selection now covers B3:B43 with B3 active,
#    and recalc so the dimension / shared formulas are all consistent.
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("B3:B43"), $true)

$excel.Calculate()
